# Applies the Review_128 -> Review_127 content replacement described by the diff.
$d = $word.ActiveDocument

# --- Paragraph 1 (Heading1): title + break + source URL -----------------
$d.Paragraphs.Item(1).Range.Text = "Review 127: [Short] Dual-Stream Diffusion Net for Text-to-Video Generation, 17.08.23" + [char]11 + "https://huggingface.co/papers/2308.08316"

# --- Paragraph 2: bold "Paper: <arxiv link>" line ------------------------
$d.Paragraphs.Item(2).Range.Text = "Paper: https://arxiv.org/abs/2308.08316v3"

# --- Paragraph 5: intro text, now only followed by a single line break --
$d.Paragraphs.Item(5).Range.Text = "הגישות הראשונות ליצירה של וידאו מטקסט באמצעות מודלי דיפוזיה יצאו לראשונה לפני שנה וחצי ומאז השתדרגו משמעותית מבחינת איכות הוידאו, אורכו והתאמתו לתיאור. היום ניתן לעשות גם פרסונליזציה למודלים אלו כלומר ליצור וידאו עם אובייקט ספציפי (החתול שלכם). " + [char]11

# --- Capture formatting templates before paragraph 6's text changes -----
$bodyTemplate = $d.Paragraphs.Item(6).Range.FormattedText
$emptyTemplate = $d.Paragraphs.Item(4).Range.FormattedText

# --- Insert the three new trailing paragraphs, seeded with formatting ---
$d.Paragraphs.Item(6).Range.InsertParagraphAfter()
$d.Paragraphs.Item(7).Range.FormattedText = $bodyTemplate

$d.Paragraphs.Item(7).Range.InsertParagraphAfter()
$d.Paragraphs.Item(8).Range.FormattedText = $bodyTemplate

$d.Paragraphs.Item(8).Range.InsertParagraphAfter()
$d.Paragraphs.Item(9).Range.FormattedText = $emptyTemplate

# --- Fill in the final text for paragraph 6 (replaced) and the two new --
# --- paragraphs; paragraph 9 stays an empty run.                       --
$d.Paragraphs.Item(6).Range.Text = "המאמר שנסקור היום ב-#shorthebrewpapereviews משדרג את הגישה הזו ומאפשר ליצור וידאו לא רק לאובייקט מסוים אלא גם לדפוס תנועה מסיום (הנגזר מוידאו אחר למשל). וכל זה בהתאמה לתיאור הטקסטואלי. איך מאמנים מודל כזה? מזינים למודל את הוידאו ובשלב הראשון מעבירים כל פריים דרך האנקודר להפקה של ייצוגו הלטנטי (עם VQ-VAE). "
$d.Paragraphs.Item(7).Range.Text = "מהייצוג הזה מאמנים מודל המפרק את הייצוג הזה את ייצוג התנועה בוידאו (בין הפריימים) לבין ייצוג התוכן של הוידאו (כל אחד מהם הוא מערך של וקטורי ייצוג) – זה נעשה באמצעות Motion Decomposer. מערך וקטורים אלו מוזן למודל דיפוזיה משלו (מכאן בא השם dual stream) שעושים את קסמיהם הרגילים. פלטי מודלי דיפוזיה אלו מוזנים לרשת ש״מערבבת״ אותם ומוציאה שני ייצוגים מסונכרנים של תנועה ושל התוכן. "
$d.Paragraphs.Item(8).Range.Text = "בסוף שני ייצוגים מסוכנרנים אלו מוזן לרשת המשלבת אותם ובונה ייצוג של וידאו שעובר דרך הדקודר כדי לגנרט וידאו. כאשר רוצים לגנרט וידאו לייצוג תנועה נתון מכיילים את המודל על ידי מזעור לוס השחזור את ייצוג התנועה מהייצוג הוידאו המגונרט."

Write-Output ("paragraphs=" + $d.Paragraphs.Count)
